$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (styles) from column P (3:33) into the new column Q
$ws.Range("P3:P33").Copy() | Out-Null
$ws.Range("Q3:Q33").PasteSpecial(-4122) | Out-Null

$values = @(
    53.463696812512026,
    46.05204738706685,
    60.76705279190513,
    46.609654277955656,
    39.785591828762811,
    53.693996785869842,
    49.132459991853935,
    42.132308166831223,
    56.225753650646354,
    28.457427087863305,
    20.524708126577082,
    36.325895173845353,
    37.816151622141014,
    29.032279844170926,
    46.928626462141906,
    51.38232216208695,
    45.862881450184311,
    57.0280888993139,
    44.951834666409091,
    38.216466887636237,
    51.83682668469686,
    82.176148450436926,
    66.965035434789911,
    96.931980629894966,
    56.391242440049062,
    50.844030930786069,
    61.300998533028128,
    54.829571415516767,
    58.407045187583961,
    51.452932817170577
)

$ws.Range("Q3").Value = 2020

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 4
    $ws.Cells.Item($row, 17).Value = $values[$i]
}

$ws.Range("S34").Select()
